# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
# Adds two new "Estado de Cuenta" detail rows (employee 29568655 DIANA CECILIA
# BONILLA BANAVIDES and a second entry for the existing employee 94372349
# JOHNNY FABIAN TORRES APARICIO) for period 2509, updates the summary
# counters (Cant. Trabajadores / Cant. Periodos) and the total "Valor Mora".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the two new detail rows right after the existing
#        data row (row 16). This pushes the trailing signature block
#        (rows 21-22) down to rows 23-24, exactly like inserting two rows.
$ws.Range("A17:A18").EntireRow.Insert()

# --- 2. Clone the formatting of the existing data row (row 16) into the
#        two freshly inserted rows so borders/fills/fonts match the table.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# --- 3. Fill in the new data.
# Row 17: new worker (29568655 / DIANA CECILIA BONILLA BANAVIDES), period 2509
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "29568655"
$ws.Range("D17").Value = "DIANA CECILIA BONILLA BANAVIDES"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: existing worker (94372349 / JOHNNY FABIAN TORRES APARICIO), period 2509
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "94372349"
$ws.Range("D18").Value = "JOHNNY FABIAN TORRES APARICIO"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# --- 4. Update the summary block.
# Cant. Trabajadores (C13) and Cant. Periodos (F13): 1 -> 2
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# Valor Mora total (E11): sum of the three detail rows' "Valor Mora" column.
$ws.Range("E11").Value = 154480
